# Update scripts with new TPM-derived values (regenerated NATMI output).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (FAPs -> ECs) ---
$ws.Range("G2").Value = 0.6811249999999999
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3251496666666667
$ws.Range("N2").Value = 0.975449
$ws.Range("O2").Value = 0.07121046526627427
$ws.Range("P2").Value = 0.07121046526627427
$ws.Range("Q2").Value = 0.2214675667083333
$ws.Range("R2").Value = 1.993208100375
$ws.Range("S2").Value = 0.07121046526627427
$ws.Range("T2").Value = 0.07121046526627427

# --- Row 3 (FAPs -> FAPs) ---
$ws.Range("G3").Value = 0.6811249999999999
$ws.Range("O3").Value = 0.2207208394324094
$ws.Range("P3").Value = 0.2207208394324094
$ws.Range("Q3").Value = 0.6864511704583331
$ws.Range("R3").Value = 6.178060534124999
$ws.Range("S3").Value = 0.2207208394324094
$ws.Range("T3").Value = 0.2207208394324094

# --- Row 4 (FAPs -> MuSCs) ---
$ws.Range("G4").Value = 0.6811249999999999
$ws.Range("M4").Value = 3.226895
$ws.Range("N4").Value = 9.680685
$ws.Range("O4").Value = 0.7067166842615477
$ws.Range("P4").Value = 0.7067166842615475
$ws.Range("Q4").Value = 2.197918856875
$ws.Range("R4").Value = 19.781269711875
$ws.Range("S4").Value = 0.7067166842615477
$ws.Range("T4").Value = 0.7067166842615475

# --- New row 5 (FAPs -> Resolving-Mac) ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6811249999999999
$ws.Range("H5").Value = 2.043375
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006173333333333333
$ws.Range("N5").Value = 0.01852
$ws.Range("O5").Value = 0.001352011039768762
$ws.Range("P5").Value = 0.001352011039768762
$ws.Range("Q5").Value = 0.004204811666666666
$ws.Range("R5").Value = 0.03784330499999999
$ws.Range("S5").Value = 0.001352011039768762
$ws.Range("T5").Value = 0.001352011039768762
